$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.909.06'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.042.18'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.90'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.667'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.16%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.59'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +8.21%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '60.94'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.42%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.385'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0789'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.59%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.69%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '16.12'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +4.32%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.337.69'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -5.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.60'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.042.49'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '36.853.23'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.76'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +15.99%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '74.71'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +2.10%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0904'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +5.86%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.74'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.92%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.86%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +10.56%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.85'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.27'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '20.16'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.64%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.72%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.14'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +2.95%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.73'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.92%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0617'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.47'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.25%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0867'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -3.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.25'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -1.04%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.77%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.111'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +14.13%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.36'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '17.79'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '96.84'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.90%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.85'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.52%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +13.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.47'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.282.38'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.89%  '
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.75'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.224.21'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.49%  '
